$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120, shifting existing rows 120-148 down to 121-149
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new data record
$ws.Cells.Item(120, 1).Value = 11
$ws.Cells.Item(120, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(120, 3).Value = "Bíobío"
$ws.Cells.Item(120, 4).Value = 44785
$ws.Cells.Item(120, 5).Value = 8
$ws.Cells.Item(120, 6).Value = 100112043
$ws.Cells.Item(120, 7).Value = "Pepino ensalada"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 100
$ws.Cells.Item(120, 11).Value = 20000
$ws.Cells.Item(120, 12).Value = 22000
$ws.Cells.Item(120, 13).Value = 21000
$ws.Cells.Item(120, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(120, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(120, 16).Value = 420
$ws.Cells.Item(120, 17).Value = 50
$ws.Cells.Item(120, 18).Value = "Hortaliza"
